# Daily attendance processing - 2026-01-13 13:56:43
# For each data row in the "Recorded By" column (G), if the value is a
# comma-separated list of 2+ names whose first entry is not "System",
# swap the first two entries (leaving any further entries untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item(1, 7).End(-4121).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Length -ge 2 -and $parts[0] -ne "System") {
            if ($parts.Length -gt 2) {
                $rest = $parts[2..($parts.Length - 1)]
                $newParts = @($parts[1], $parts[0]) + $rest
            } else {
                $newParts = @($parts[1], $parts[0])
            }
            $newVal = $newParts -join ", "
            $cell.Value = $newVal
        }
    }
}
